$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @("Veículo", "Frota", "Momento Infração", "Infração", "criticidade", "Evidência", "Local Infração", "Situação", "Status")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Copy the header formatting (bold, border, centered) from A1 onto the
# newly-added header cells E1:I1 so they match the rest of the row.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1:I1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 2
$ws.Cells.Item(2, 1).Value = "52001 - QYN4G48"
$ws.Cells.Item(2, 2).Value = "RCR-PE"
$ws.Cells.Item(2, 3).Value = "7 de fev. de 2025, 14:42:16"
$ws.Cells.Item(2, 4).Value = "Motorista com celular"
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = "https://moovsec-videos-prod.s3.sa-east-1.amazonaws.com/video_evidence_67a6462d2f869d1ad65787e0.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Date=20250208T070005Z&X-Amz-SignedHeaders=host&X-Amz-Expires=604800&X-Amz-Credential=AKIA3QTD5B6Z3DVDF6WA%2F20250208%2Fsa-east-1%2Fs3%2Faws4_request&X-Amz-Signature=4650e66ae32657b678b990ca908fc6cc27a1899625536a336266aae9703eeebd"
$ws.Cells.Item(2, 7).Value = "https://www.google.com/maps?q=-7.56624600000000700,-34.99236800000000000"
$ws.Cells.Item(2, 8).Value = "Verdadeiro"
$ws.Cells.Item(2, 9).Value = "Validado"

# Row 3
$ws.Cells.Item(3, 1).Value = "322238 - QZX5G06"
$ws.Cells.Item(3, 2).Value = "RCR-AM"
$ws.Cells.Item(3, 3).Value = "5 de fev. de 2025, 08:28:55"
$ws.Cells.Item(3, 4).Value = "Motorista com celular"
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = "https://moovsec-videos-prod.s3.sa-east-1.amazonaws.com/video_evidence_67a34b8ebcc233ee604f3ed1.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Date=20250206T063801Z&X-Amz-SignedHeaders=host&X-Amz-Expires=604800&X-Amz-Credential=AKIA3QTD5B6Z3DVDF6WA%2F20250206%2Fsa-east-1%2Fs3%2Faws4_request&X-Amz-Signature=c70eaa1b3fcee66d541bd27f6e1a56a1e032fbfdeb5f758d537d82d8a1c5ea8b"
$ws.Cells.Item(3, 7).Value = "https://www.google.com/maps?q=-3.05631099999999400,-59.98609100000000000"
$ws.Cells.Item(3, 8).Value = "Erro"
$ws.Cells.Item(3, 9).Value = "Erro"

# Row 4
$ws.Cells.Item(4, 1).Value = "322232 - QZX5E76"
$ws.Cells.Item(4, 2).Value = "RCR-AM"
$ws.Cells.Item(4, 3).Value = "8 de fev. de 2025, 07:37:51"
$ws.Cells.Item(4, 4).Value = "Motorista com celular"
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = "https://moovsec-videos-prod.s3.sa-east-1.amazonaws.com/video_evidence_67a7343e2f869d1ad68480e2.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Date=20250209T062314Z&X-Amz-SignedHeaders=host&X-Amz-Expires=604799&X-Amz-Credential=AKIA3QTD5B6Z3DVDF6WA%2F20250209%2Fsa-east-1%2Fs3%2Faws4_request&X-Amz-Signature=9f6aea201adeab09ea216546fbe7657eb15e9c6b935d75ad04b9a60fdf70601f"
$ws.Cells.Item(4, 7).Value = "https://www.google.com/maps?q=-3.04053000000000400,-59.92979300000000400"
$ws.Cells.Item(4, 8).Value = "Verdadeiro"
$ws.Cells.Item(4, 9).Value = "Validado"
